# Automatische test-sync: 2025-06-29 13:59:50
$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new mail-log entry as row 4 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(4, 1).Value = "Wanneer zijn jullie open?"
$logs.Cells.Item(4, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(4, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Cells.Item(4, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item(4, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item(4, 6).Value = "2025-06-29 13:59:27"
$logs.Cells.Item(4, 7).Value = "Ja"
$logs.Cells.Item(4, 8).Value = "Nee"
$logs.Cells.Item(4, 9).Value = "Ja"

# --- extend the conditional-formatting ranges to cover the new row ---
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range("$col`2:$col`3")
    $newRange = $logs.Range("$col`2:$col`4")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the "Openingstijden / Locatie" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 2
